# Plantilla de casos de uso completada
# Fills in the previously-empty "Descripcion de Caso de Uso" (column C) cells
# for CU-02, CU-04, CU-05, CU-08, CU-09, CU-11 and CU-16 on the "Casos de Uso"
# sheet, and restores the view/selection state left by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

$ws.Range("C6").Value  = 'El director debe poder manipular la información de los colaboradores en el sistema, para así poder relacionar pagos, asignar grupos y horarios de las actividades respectivas a cada colaborador'
$ws.Range("C8").Value  = 'El director debe poder  visualizar, registrar y editar los datos de rentas de aulas, con el fin de llevar un control de los espacios y su disponibilidad'
$ws.Range("C9").Value  = 'El director debe poder  asignar lapsos a una actividad a ser realizada, con el fin de evitar traslapes y conocer en todo momento la disponibilidad de los espacios'
$ws.Range("C12").Value = 'El director debe poder almacenar los pagos recibidos por parte de los alumnos para sus colaboradores para ampliar la disponibilidad de los cobros'
$ws.Range("C13").Value = 'El director debe poder visualizar los detalles de los pagos que ha recibido por parte de los alumnos para saber con certeza los movimientos pendientes de realizar'
$ws.Range("C15").Value = 'El director debe poder almacenar los datos de las campañas publicitarias que se han creado, con el fin de tener actualizados en todo momento la duración, campañas activas y egresos relacionados a las campañas'
$ws.Range("C20").Value = 'El maestro debe poder registrar los pagos realizados por los alumnos, para llevar un control de ingresos y pagos pendientes'

# Rows 6 and 15 now wrap onto a third line in the wide column C, so Excel's
# autofit grows their height from 30pt to 45pt.
$ws.Rows("6").RowHeight = 45
$ws.Rows("15").RowHeight = 45

# Restore the scroll position / active selection left in the saved file.
$ws.Activate()
$ws.Range("C20").Select()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 3

$wb.Save()
